$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- The edit: swap the two task rows (row 3 <-> row 4), i.e. the task that used
# to be "Prescribe Medication" (id 2) and "herbie_basic_test" (id 3) trade places
# in the list, while keeping the id column (A) fixed. Columns B:E (name,
# description, url, script) plus their cell styles, the row heights, and the
# hyperlinks on column D all move together with the row's content.

# 1) Detach existing hyperlinks before we start shuffling cells around so that
#    re-creating them afterwards doesn't collide with stale relationship ids.
#    (Deleting the hyperlinks on one cell clears the sheet's hyperlink list.)
$ws.Range("D3").Hyperlinks.Delete()

# 2) Re-create the two hyperlinks already pointing at the *destination* cell for
#    each URL, so the link target correctly follows its text to the new row.
#    D3 held the "pat_id=18" link before the swap -> that link now belongs on D4.
#    D4 held the "login.html" link before the swap -> that link now belongs on D3.
$ws.Hyperlinks.Add($ws.Range("D4"), "https://hrithik.webchartnow.com/webchart.cgi?f=chart&s=pat&pat_id=18")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://mieweb.github.io/herbie/playgrounds/login.html")

# 3) Swap the B:E cell contents (values + number formats/styles) of row 3 and
#    row 4 using a scratch range as a temporary holding area. This also
#    re-stamps D3/D4 with their original per-cell styles, overwriting whatever
#    formatting Hyperlinks.Add applied in step 2.
$ws.Range("B3:E3").Copy($ws.Range("B100:E100"))
$ws.Range("B4:E4").Copy($ws.Range("B3:E3"))
$ws.Range("B100:E100").Copy($ws.Range("B4:E4"))
$ws.Range("B100:E100").Clear()

# 4) Swap the row heights so they travel with the content that moved.
$height3 = $ws.Rows("3:3").RowHeight
$height4 = $ws.Rows("4:4").RowHeight
$ws.Rows("3:3").RowHeight = $height4
$ws.Rows("4:4").RowHeight = $height3

# 5) Leave the selection where the author's editing session ended up.
$ws.Range("I4").Select()
